# Daily rolling update of the BP terminal gate pricing table.
# The newest effective date advances from 20 Feb 2026 (46073) to 21 Feb 2026 (46074),
# each existing date block shifts down to the next-older date, and new prices are set.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("A8").Value = 46074
$ws.Range("D8").Value = 158.09
$ws.Range("E8").Value = 147.71
$ws.Range("F8").Value = 157.71
$ws.Range("G8").Value = 147.6

# Row 9
$ws.Range("A9").Value = 46074
$ws.Range("D9").Value = 158.09
$ws.Range("E9").Value = 147.71
$ws.Range("F9").Value = 157.71
$ws.Range("G9").Value = 147.6

# Row 10
$ws.Range("A10").Value = 46074
$ws.Range("D10").Value = 159.72
$ws.Range("E10").Value = 150.41
$ws.Range("F10").Value = 160.41
$ws.Range("G10").Value = 150.65

# Row 11
$ws.Range("A11").Value = 46073
$ws.Range("D11").Value = 157.41
$ws.Range("E11").Value = 147.24
$ws.Range("F11").Value = 157.25
$ws.Range("G11").Value = 147.13

# Row 12
$ws.Range("A12").Value = 46073
$ws.Range("D12").Value = 157.41
$ws.Range("E12").Value = 147.24
$ws.Range("F12").Value = 157.25
$ws.Range("G12").Value = 147.13

# Row 13
$ws.Range("A13").Value = 46073
$ws.Range("D13").Value = 158.87
$ws.Range("E13").Value = 149.98
$ws.Range("F13").Value = 159.98
$ws.Range("G13").Value = 150.22

# Row 17
$ws.Range("A17").Value = 46074
$ws.Range("D17").Value = 164.13
$ws.Range("E17").Value = 153.99
$ws.Range("F17").Value = 163.99

# Row 18
$ws.Range("A18").Value = 46073
$ws.Range("D18").Value = 163.27
$ws.Range("E18").Value = 153.55
$ws.Range("F18").Value = 163.55

# Row 22
$ws.Range("A22").Value = 46074
$ws.Range("D22").Value = 159.49
$ws.Range("E22").Value = 150.06
$ws.Range("F22").Value = 159.66
$ws.Range("G22").Value = 151.81

# Row 23
$ws.Range("A23").Value = 46074
$ws.Range("D23").Value = 164.69
$ws.Range("E23").Value = 156.11
$ws.Range("F23").Value = 166.11

# Row 24
$ws.Range("A24").Value = 46074
$ws.Range("D24").Value = 164.88
$ws.Range("E24").Value = 156.63
$ws.Range("F24").Value = 166.63

# Row 25
$ws.Range("A25").Value = 46074
$ws.Range("D25").Value = 164.88
$ws.Range("E25").Value = 156.14
$ws.Range("F25").Value = 166.14
$ws.Range("G25").Value = 157

# Row 26
$ws.Range("A26").Value = 46074
$ws.Range("D26").Value = 164.52
$ws.Range("E26").Value = 157.72
$ws.Range("F26").Value = 167.72

# Row 27
$ws.Range("A27").Value = 46073
$ws.Range("D27").Value = 158.59
$ws.Range("E27").Value = 149.59
$ws.Range("F27").Value = 159.19
$ws.Range("G27").Value = 151.34

# Row 28
$ws.Range("A28").Value = 46073
$ws.Range("D28").Value = 163.83
$ws.Range("E28").Value = 155.68
$ws.Range("F28").Value = 165.68

# Row 29
$ws.Range("A29").Value = 46073
$ws.Range("D29").Value = 164.03
$ws.Range("E29").Value = 156.21
$ws.Range("F29").Value = 166.21

# Row 30
$ws.Range("A30").Value = 46073
$ws.Range("D30").Value = 164.03
$ws.Range("E30").Value = 155.71
$ws.Range("F30").Value = 165.71
$ws.Range("G30").Value = 156.57

# Row 31
$ws.Range("A31").Value = 46073
$ws.Range("D31").Value = 163.67
$ws.Range("E31").Value = 157.29
$ws.Range("F31").Value = 167.29

# Row 35
$ws.Range("A35").Value = 46074
$ws.Range("D35").Value = 158.11
$ws.Range("E35").Value = 148.11
$ws.Range("F35").Value = 157.11

# Row 36
$ws.Range("A36").Value = 46073
$ws.Range("D36").Value = 157.37
$ws.Range("E36").Value = 147.68
$ws.Range("F36").Value = 156.68

# Row 40
$ws.Range("A40").Value = 46074
$ws.Range("D40").Value = 164.39
$ws.Range("E40").Value = 155.28
$ws.Range("F40").Value = 165.28

# Row 41
$ws.Range("A41").Value = 46074
$ws.Range("D41").Value = 164.11
$ws.Range("E41").Value = 155.7
$ws.Range("F41").Value = 165.7

# Row 42
$ws.Range("A42").Value = 46073
$ws.Range("D42").Value = 163.54
$ws.Range("E42").Value = 154.75
$ws.Range("F42").Value = 164.75

# Row 43
$ws.Range("A43").Value = 46073
$ws.Range("D43").Value = 163.26
$ws.Range("E43").Value = 155.17
$ws.Range("F43").Value = 165.17

# Row 47
$ws.Range("A47").Value = 46074
$ws.Range("D47").Value = 158.59
$ws.Range("E47").Value = 150.31
$ws.Range("F47").Value = 160.31

# Row 48
$ws.Range("A48").Value = 46074
$ws.Range("D48").Value = 158.3
$ws.Range("E48").Value = 150.3
$ws.Range("F48").Value = 160.3

# Row 49
$ws.Range("A49").Value = 46073
$ws.Range("D49").Value = 158.43
$ws.Range("E49").Value = 150.37
$ws.Range("F49").Value = 160.37

# Row 50
$ws.Range("A50").Value = 46073
$ws.Range("D50").Value = 158.14
$ws.Range("E50").Value = 150.36
$ws.Range("F50").Value = 160.36

# Row 54
$ws.Range("A54").Value = 46074
$ws.Range("D54").Value = 173.62
$ws.Range("E54").Value = 162.86
$ws.Range("F54").Value = 172.86

# Row 55
$ws.Range("A55").Value = 46074
$ws.Range("D55").Value = 162.97
$ws.Range("E55").Value = 162.04
$ws.Range("F55").Value = 172.04

# Row 56
$ws.Range("A56").Value = 46074
$ws.Range("D56").Value = 162.61

# Row 57
$ws.Range("A57").Value = 46074
$ws.Range("D57").Value = 163.55
$ws.Range("E57").Value = 156.46

# Row 58
$ws.Range("A58").Value = 46074
$ws.Range("D58").Value = 159.32
$ws.Range("E58").Value = 152.36
$ws.Range("F58").Value = 162.36

# Row 59
$ws.Range("A59").Value = 46074
$ws.Range("D59").Value = 166.57
$ws.Range("E59").Value = 161.46

# Row 60
$ws.Range("A60").Value = 46073
$ws.Range("D60").Value = 172.75
$ws.Range("E60").Value = 162.47
$ws.Range("F60").Value = 172.47

# Row 61
$ws.Range("A61").Value = 46073
$ws.Range("D61").Value = 162.11
$ws.Range("E61").Value = 161.71
$ws.Range("F61").Value = 171.71

# Row 62
$ws.Range("A62").Value = 46073
$ws.Range("D62").Value = 161.98

# Row 63
$ws.Range("A63").Value = 46073
$ws.Range("D63").Value = 162.91
$ws.Range("E63").Value = 156.13

# Row 64
$ws.Range("A64").Value = 46073
$ws.Range("D64").Value = 158.68
$ws.Range("E64").Value = 152.03
$ws.Range("F64").Value = 162.03

# Row 65
$ws.Range("A65").Value = 46073
$ws.Range("D65").Value = 165.71
$ws.Range("E65").Value = 161.05
